$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.163.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +11.04%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.270.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +6.07%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "397.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.562"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.91%  "

$ws.Range("E8").Value = "  -0.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.627"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.38"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0962"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +12.67%  "

$ws.Range("E12").Value = "  +2.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.775.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.76%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.260.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.05%  "

$ws.Range("E17").Value = "  +1.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "56.844.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +10.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000106"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.86%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "308.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +16.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.46%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.25%  "

$ws.Range("E27").Value = "  +5.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.07%  "

$ws.Range("E29").Value = "  +2.84%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.92%  "

$ws.Range("E31").Value = "  -0.36%  "

$ws.Range("E32").Value = "  +5.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "37.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0482"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.66%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.56"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.44%  "

$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +22.39%  "

$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.52"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.06%  "

$ws.Range("E40").Value = "  -0.08%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "134.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.07%  "

$ws.Range("E43").Value = "  -0.66%  "

$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.98%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.120"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.282"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.69%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.151.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.81%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.84%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.37"
$ws.Range("D50").Style = "Normal"

$ws.Range("E51").Value = "  +37.53%  "
